$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 2.336779832839966
$ws.Range("B1").Value = 5.437224388122559
$ws.Range("C1").Value = 2.32887601852417
$ws.Range("D1").Value = 1.602648854255676
$ws.Range("E1").Value = 1.464346051216125
